$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.3398761212379742
$ws1.Range("C2").Value = -0.2080425783312123
$ws1.Range("B3").Value = -1.178633112448592
$ws1.Range("C3").Value = -0.4029934313610458
$ws1.Range("B4").Value = -0.7418062170048667
$ws1.Range("C4").Value = -0.05510495897711093

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.7652778581602446
$ws2.Range("C2").Value = -0.05331725020813981
$ws2.Range("B3").Value = -1.108795268790843
$ws2.Range("C3").Value = -0.506664820291607
$ws2.Range("B4").Value = -0.2154940908642698
$ws2.Range("C4").Value = 0.4137171347917183
